$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column H header "kierunek" and fill rows 2-6 with "matematyka stosowana"
$ws.Range("H1").Value = "kierunek"
$ws.Range("H2").Value = "matematyka stosowana"
$ws.Range("H3").Value = "matematyka stosowana"
$ws.Range("H4").Value = "matematyka stosowana"
$ws.Range("H5").Value = "matematyka stosowana"
$ws.Range("H6").Value = "matematyka stosowana"

# Copy style of G1 (header style) to H1, then add fill to match new style
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats

# Select cell to match diff's reported selection
$ws.Range("I7").Select()
